$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.424.37'
$ws.Range("E2").Value = '  -0.51%  '

$ws.Range("D3").Value = '2.285.04'
$ws.Range("E3").Value = '  -0.06%  '

$ws.Range("E4").Value = '  -0.33%  '

$ws.Range("D5").Value = '113.78'
$ws.Range("E5").Value = '  +1.67%  '

$ws.Range("D6").Value = '266.84'
$ws.Range("E6").Value = '  -0.38%  '

$ws.Range("D7").Value = '0.621'
$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("E9").Value = '  -0.34%  '

$ws.Range("D10").Value = '47.88'
$ws.Range("E10").Value = '  +1.60%  '

$ws.Range("D11").Value = '0.0933'
$ws.Range("E11").Value = '  -0.17%  '

$ws.Range("D12").Value = '9.37'
$ws.Range("E12").Value = '  +11.13%  '

$ws.Range("E13").Value = '  +1.21%  '

$ws.Range("D14").Value = '15.53'
$ws.Range("E14").Value = '  +0.27%  '

$ws.Range("D15").Value = '2.614.27'
$ws.Range("E15").Value = '  -0.66%  '

$ws.Range("D16").Value = '0.870'
$ws.Range("E16").Value = '  +3.24%  '

$ws.Range("D17").Value = '2.283.85'
$ws.Range("E17").Value = '  -0.58%  '

$ws.Range("D18").Value = '43.313.99'
$ws.Range("E18").Value = '  -0.56%  '

$ws.Range("E19").Value = '  -0.53%  '

$ws.Range("D20").Value = '6.88'
$ws.Range("E20").Value = '  +5.39%  '

$ws.Range("E21").Value = '  -0.67%  '

$ws.Range("E22").Value = '  +1.07%  '

$ws.Range("D23").Value = '234.20'
$ws.Range("E23").Value = '  +0.93%  '

$ws.Range("E24").Value = '  +3.62%  '

$ws.Range("D25").Value = '2.93'
$ws.Range("E25").Value = '  +3.88%  '

$ws.Range("D27").Value = '11.44'
$ws.Range("E27").Value = '  +0.90%  '

$ws.Range("D28").Value = '3.99'
$ws.Range("E28").Value = '  +0.53%  '

$ws.Range("D29").Value = '40.88'
$ws.Range("E29").Value = '  -3.58%  '

$ws.Range("E30").Value = '  -2.54%  '

$ws.Range("E31").Value = '  -0.60%  '

$ws.Range("D32").Value = '173.38'
$ws.Range("E32").Value = '  -1.53%  '

$ws.Range("D33").Value = '21.50'
$ws.Range("E33").Value = '  -0.32%  '

$ws.Range("D34").Value = '0.0907'
$ws.Range("E34").Value = '  -1.13%  '

$ws.Range("D35").Value = '5.75'
$ws.Range("E35").Value = '  +5.26%  '

$ws.Range("E36").Value = '  +0.94%  '

$ws.Range("E37").Value = '  -0.13%  '

$ws.Range("D38").Value = '0.0369'
$ws.Range("E38").Value = '  +4.65%  '

$ws.Range("D39").Value = '3.94'
$ws.Range("E39").Value = '  +4.88%  '

$ws.Range("E40").Value = '  -3.56%  '

$ws.Range("D41").Value = '2.70'
$ws.Range("E41").Value = '  +12.49%  '

$ws.Range("D42").Value = '78.06'
$ws.Range("E42").Value = '  +7.50%  '

$ws.Range("D43").Value = '14.32'
$ws.Range("E43").Value = '  +6.01%  '

$ws.Range("D44").Value = '0.240'
$ws.Range("E44").Value = '  +0.20%  '

$ws.Range("E45").Value = '  +6.51%  '

$ws.Range("E46").Value = '  -0.25%  '

$ws.Range("E47").Value = '  -1.64%  '

$ws.Range("D48").Value = '8.72'
$ws.Range("E48").Value = '  +0.25%  '

$ws.Range("D49").Value = '104.43'
$ws.Range("E49").Value = '  +2.06%  '

$ws.Range("E50").Value = '  +2.97%  '

$ws.Range("E51").Value = '  +0.11%  '
